$wb = $excel.ActiveWorkbook

# The original sheet ("Test") carries legacy "preserved shared-formula group"
# metadata from import, which keeps re-emitting <f t="shared" .../> on save
# even after the formulas are rewritten. Replace the sheet wholesale with a
# fresh one so formulas serialize as plain, independent <f> elements (as in
# the target workbook) and the old revision/"xr" scaffolding is dropped.
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "TempNewSheet"
$wb.Worksheets.Item("Test").Delete()
$wsNew.Name = "Sheet1"

$wsNew.Range("A1").Value = 0
$wsNew.Range("B1").Formula = "=MOD(A1,2)=0"
$wsNew.Range("C1").Formula = "=A2/A1"

$wsNew.Range("A2").Value = 1
$wsNew.Range("B2").Formula = "=MOD(A2,2)=0"
$wsNew.Range("C2").Formula = "=A3/A2"

$wsNew.Range("A3").Value = 2
$wsNew.Range("B3").Formula = "=MOD(A3,2)=0"
$wsNew.Range("C3").Formula = "=A4/A3"

$wsNew.Range("A4").Value = 3
$wsNew.Range("B4").Formula = "=MOD(A4,2)=0"
$wsNew.Range("C4").Formula = "=A5/A4"

$wsNew.Range("A5").Formula = "=SUM(A1:A4)"
$wsNew.Range("B5").Formula = "=MOD(A5,2)=0"
$wsNew.Range("C5").Formula = "=A6/A5"

$wsNew.Range("A6").Value = "a"
$wsNew.Range("B6").Value = "b"
$wsNew.Range("C6").Formula = "=A6&B6"

# Rows 7-9 (external-link / #NAME? / #NULL! demo formulas) and the
# external-workbook link they (partly) relied on are removed entirely.
$wb.BreakLink("OTHER", 1)

# Restore the active selection shown in the saved workbook.
$wsNew.Range("G11").Select()
